$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New unstyled "Address" datatype row (mirrors the plain style of rows 5-6)
$ws.Range("B7").Value = "Address"
$ws.Range("C7").Value = "adr"

# New "Environment" / "import" block (rows 10-11), formatted like the bordered
# header block (B1:D4) so the cells carry an explicit style instead of the
# workbook default.
$ws.Range("B10:C11").NumberFormat = "General"
$ws.Range("B10").Value = "Environment"
$ws.Range("B11").Value = "import"
$ws.Range("C11").Value = "com.example.beans"

# Column widths roughly matching the autosized (best-fit) widths from the
# target sheet (B=15.71, C=19.14, D=25 "characters")
$ws.Range("B:B").ColumnWidth = 14.8
$ws.Range("C:C").ColumnWidth = 18.35
$ws.Range("D:D").ColumnWidth = 24.17

# Selection / active cell as recorded in the saved view
$ws.Range("C10").Select()

$wb.Save()
